$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Update the wording of a few cells in row 3 (they now represent a
# value that is shared/merged across rows 3 and 4 for columns C, D, F)
# ------------------------------------------------------------------
$ws.Range("C3").Value = "data2,3-2"
$ws.Range("D3").Value = "data2,3-3"
$ws.Range("F3").Value = "data2,3-5"

# The corresponding cells in row 4 are no longer needed once the
# columns are merged vertically with row 3, so clear their contents.
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("F4").ClearContents()

# ------------------------------------------------------------------
# Apply left alignment to the merged columns, and highlight column C
# with a yellow fill.
# ------------------------------------------------------------------
$ws.Range("D3:D4").HorizontalAlignment = -4131
$ws.Range("F3:F4").HorizontalAlignment = -4131

$ws.Range("C3:C4").HorizontalAlignment = -4131
$ws.Range("C3:C4").Interior.ColorIndex = 6

# ------------------------------------------------------------------
# Merge the cells so the single value spans both rows.
# ------------------------------------------------------------------
$ws.Range("C3:C4").Merge()
$ws.Range("D3:D4").Merge()
$ws.Range("F3:F4").Merge()

# ------------------------------------------------------------------
# Clear the lingering selection rectangle that used to highlight B4.
# ------------------------------------------------------------------
$ws.Range("A1").Select()
